# Automatic update of files.
#
# The underlying observation records for rows 9-12 (one cluster of
# "Kyrkberget, Dlr" sightings) and rows 23,24,26-30 (a second cluster) were
# re-sequenced: the species/occurrence data landed on different row numbers
# than before, while everything else about each row (location columns,
# dates, observer, formatting, the blank "marker" cells, etc.) stayed put.
#
# Concretely, column set A,B,D,E,F,G,H,Q,R (Id, Taxonsorteringsordning,
# Rödlistade, TaxonId, Artnamn, Vetenskapligt namn, Auktor, Ost, Nord) is
# permuted across the two row clusters; nothing else in the sheet changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The columns whose values travel together as a single observation record.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Destination row -> source row: the record currently sitting in the source
# row (as the workbook stands right now, before this script runs) is the one
# that ends up in the destination row.
$rowMap = @{
    9  = 12
    10 = 11
    11 = 10
    12 = 9
    23 = 27
    24 = 29
    26 = 23
    27 = 24
    28 = 26
    29 = 30
    30 = 28
}

# 1) Snapshot every involved row's current values first - this is a cyclic
#    permutation (rows both donate and receive data), so nothing may be
#    written until all the old values have been read.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each destination row's record from the snapshot taken above.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}

# 3) The empty "Bestämningsmetod" placeholder cell (column AF) that used to
#    sit on row 12 travels with that row's record, which now lives on row 9.
$ws.Range("AF12").ClearContents()
$ws.Range("AF9").NumberFormat = "General"

Write-Output "Row records permuted across clusters [9-12] and [23,24,26-30]."
